# Update cryptocurrency price / volume figures (and restore the swapped
# Fetch.AI / WEMIXToken rows) to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Formula = '="69.385.41"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Formula = '="  +1.64%  "'
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)

# Row 3
$ws.Range("D3").Formula = '="3.932.88"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Formula = '="  +0.30%  "'
$ws.Range("E3").Copy()
$ws.Range("E3").PasteSpecial(-4163)

# Row 4
$ws.Range("E4").Formula = '="  +0.05%  "'
$ws.Range("E4").Copy()
$ws.Range("E4").PasteSpecial(-4163)

# Row 5
$ws.Range("D5").Formula = '="515.59"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Formula = '="  +5.82%  "'
$ws.Range("E5").Copy()
$ws.Range("E5").PasteSpecial(-4163)

# Row 6
$ws.Range("D6").Formula = '="145.87"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Formula = '="  -1.17%  "'
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)

# Row 8
$ws.Range("E8").Formula = '="  -0.04%  "'
$ws.Range("E8").Copy()
$ws.Range("E8").PasteSpecial(-4163)

# Row 9
$ws.Range("D9").Formula = '="0.730"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Formula = '="  -0.57%  "'
$ws.Range("E9").Copy()
$ws.Range("E9").PasteSpecial(-4163)

# Row 10
$ws.Range("D10").Formula = '="0.172"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Formula = '="  +3.34%  "'
$ws.Range("E10").Copy()
$ws.Range("E10").PasteSpecial(-4163)

# Row 11
$ws.Range("D11").Formula = '="0.0000342"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Formula = '="  -2.11%  "'
$ws.Range("E11").Copy()
$ws.Range("E11").PasteSpecial(-4163)

# Row 12
$ws.Range("D12").Formula = '="43.07"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Formula = '="  -0.04%  "'
$ws.Range("E12").Copy()
$ws.Range("E12").PasteSpecial(-4163)

# Row 13
$ws.Range("D13").Formula = '="4.562.34"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Formula = '="  +0.22%  "'
$ws.Range("E13").Copy()
$ws.Range("E13").PasteSpecial(-4163)

# Row 14
$ws.Range("E14").Formula = '="  -3.57%  "'
$ws.Range("E14").Copy()
$ws.Range("E14").PasteSpecial(-4163)

# Row 15
$ws.Range("D15").Formula = '="3.933.27"'
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Formula = '="  +0.09%  "'
$ws.Range("E15").Copy()
$ws.Range("E15").PasteSpecial(-4163)

# Row 16
$ws.Range("E16").Formula = '="  -2.32%  "'
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)

# Row 17
$ws.Range("E17").Formula = '="  -0.64%  "'
$ws.Range("E17").Copy()
$ws.Range("E17").PasteSpecial(-4163)

# Row 18
$ws.Range("E18").Formula = '="  +7.45%  "'
$ws.Range("E18").Copy()
$ws.Range("E18").PasteSpecial(-4163)

# Row 19
$ws.Range("D19").Formula = '="19.75"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Formula = '="  -0.99%  "'
$ws.Range("E19").Copy()
$ws.Range("E19").PasteSpecial(-4163)

# Row 20
$ws.Range("D20").Formula = '="69.321.30"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Formula = '="  +1.33%  "'
$ws.Range("E20").Copy()
$ws.Range("E20").PasteSpecial(-4163)

# Row 21
$ws.Range("D21").Formula = '="432.02"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Formula = '="  -2.27%  "'
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial(-4163)

# Row 22
$ws.Range("D22").Formula = '="3.40"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Formula = '="  -2.72%  "'
$ws.Range("E22").Copy()
$ws.Range("E22").PasteSpecial(-4163)

# Row 23
$ws.Range("D23").Formula = '="14.51"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Formula = '="  -4.66%  "'
$ws.Range("E23").Copy()
$ws.Range("E23").PasteSpecial(-4163)

# Row 24
$ws.Range("D24").Formula = '="88.20"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Formula = '="  +0.02%  "'
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial(-4163)

# Row 25
$ws.Range("D25").Formula = '="11.83"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Formula = '="  +3.06%  "'
$ws.Range("E25").Copy()
$ws.Range("E25").PasteSpecial(-4163)

# Row 26
$ws.Range("D26").Formula = '="3.90"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Formula = '="  +6.90%  "'
$ws.Range("E26").Copy()
$ws.Range("E26").PasteSpecial(-4163)

# Row 27
$ws.Range("D27").Formula = '="11.10"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Formula = '="  -3.71%  "'
$ws.Range("E27").Copy()
$ws.Range("E27").PasteSpecial(-4163)

# Row 28
$ws.Range("D28").Formula = '="36.72"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Formula = '="  -4.98%  "'
$ws.Range("E28").Copy()
$ws.Range("E28").PasteSpecial(-4163)

# Row 29
$ws.Range("E29").Formula = '="  -1.34%  "'
$ws.Range("E29").Copy()
$ws.Range("E29").PasteSpecial(-4163)

# Row 30
$ws.Range("D30").Formula = '="702.07"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Formula = '="  -2.04%  "'
$ws.Range("E30").Copy()
$ws.Range("E30").PasteSpecial(-4163)

# Row 31
$ws.Range("D31").Formula = '="13.28"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Formula = '="  -4.12%  "'
$ws.Range("E31").Copy()
$ws.Range("E31").PasteSpecial(-4163)

# Row 32
$ws.Range("D32").Formula = '="0.127"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Formula = '="  -2.82%  "'
$ws.Range("E32").Copy()
$ws.Range("E32").PasteSpecial(-4163)

# Row 34
$ws.Range("D34").Formula = '="66.53"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Formula = '="  +8.11%  "'
$ws.Range("E34").Copy()
$ws.Range("E34").PasteSpecial(-4163)

# Row 35
$ws.Range("D35").Formula = '="0.442"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Formula = '="  +5.97%  "'
$ws.Range("E35").Copy()
$ws.Range("E35").PasteSpecial(-4163)

# Row 36
$ws.Range("D36").Formula = '="0.0₃0878"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Formula = '="  -0.34%  "'
$ws.Range("E36").Copy()
$ws.Range("E36").PasteSpecial(-4163)

# Row 37
$ws.Range("D37").Formula = '="5.95"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Formula = '="  -6.25%  "'
$ws.Range("E37").Copy()
$ws.Range("E37").PasteSpecial(-4163)

# Row 38
$ws.Range("D38").Formula = '="40.27"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Formula = '="  -4.77%  "'
$ws.Range("E38").Copy()
$ws.Range("E38").PasteSpecial(-4163)

# Row 39
$ws.Range("D39").Formula = '="0.148"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Formula = '="  -0.59%  "'
$ws.Range("E39").Copy()
$ws.Range("E39").PasteSpecial(-4163)

# Row 40
$ws.Range("D40").Formula = '="0.997"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Formula = '="  -0.38%  "'
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial(-4163)

# Row 41
$ws.Range("E41").Formula = '="  -0.04%  "'
$ws.Range("E41").Copy()
$ws.Range("E41").PasteSpecial(-4163)

# Row 42
$ws.Range("D42").Formula = '="0.0484"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Formula = '="  +0.59%  "'
$ws.Range("E42").Copy()
$ws.Range("E42").PasteSpecial(-4163)

# Row 43
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").Formula = '="2.82"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Formula = '="  -8.25%  "'
$ws.Range("E43").Copy()
$ws.Range("E43").PasteSpecial(-4163)

# Row 44
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").Formula = '="3.10"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Formula = '="  +6.09%  "'
$ws.Range("E44").Copy()
$ws.Range("E44").PasteSpecial(-4163)

# Row 45
$ws.Range("D45").Formula = '="3.02"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Formula = '="  -8.45%  "'
$ws.Range("E45").Copy()
$ws.Range("E45").PasteSpecial(-4163)

# Row 46
$ws.Range("E46").Formula = '="  +0.55%  "'
$ws.Range("E46").Copy()
$ws.Range("E46").PasteSpecial(-4163)

# Row 47
$ws.Range("D47").Formula = '="3.35"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Formula = '="  +0.71%  "'
$ws.Range("E47").Copy()
$ws.Range("E47").PasteSpecial(-4163)

# Row 48
$ws.Range("D48").Formula = '="0.0₆0359"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Formula = '="  +2.26%  "'
$ws.Range("E48").Copy()
$ws.Range("E48").PasteSpecial(-4163)

# Row 49
$ws.Range("E49").Formula = '="  -2.26%  "'
$ws.Range("E49").Copy()
$ws.Range("E49").PasteSpecial(-4163)

# Row 50
$ws.Range("D50").Formula = '="2.95"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Formula = '="  +3.77%  "'
$ws.Range("E50").Copy()
$ws.Range("E50").PasteSpecial(-4163)

# Row 51
$ws.Range("D51").Formula = '="2.09"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Formula = '="  -2.39%  "'
$ws.Range("E51").Copy()
$ws.Range("E51").PasteSpecial(-4163)
